$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold, border, centered) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-49
$iValues = @(9,9,6,8,9,9,7,8,8,8,8,9,9,8,8,9,7,8,9,8,9,9,8,8,8,6,9,9,8,8,8,9,9,9,6,8,9,9,7,7,8,8,4,7,6,5,4,4)
$jValues = @(9,9,6,8,9,9,8,8,8,8,8,9,9,8,8,9,7,8,9,8,9,9,9,8,8,6,9,9,8,9,9,9,9,9,8,8,10,9,7,7,8,8,4,7,6,5,4,4)

for ($r = 2; $r -le 49; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
